$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Fitness column (C) values for rows 2 through 136 to 7293,
# matching the updated run log data.
$ws.Range("C2:C136").Value = 7293
